$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coerce the Price/Volume columns to Text before writing so Excel
# does not reinterpret dotted-thousands price strings (e.g. "26.667.40")
# or leading-zero / percent strings as numbers.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '26.667.40'
$ws.Range('D3').Value = '1.600.11'
$ws.Range('E3').Value = '  +0.54%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '211.42'
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('D6').Value = '0.518'
$ws.Range('E6').Value = '  +1.21%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('E9').Value = '  -1.09%  '
$ws.Range('D10').Value = '19.52'
$ws.Range('E10').Value = '  -0.76%  '
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('D12').Value = '1.824.19'
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('D13').Value = '1.581.74'
$ws.Range('E13').Value = '  -0.76%  '
$ws.Range('D14').Value = '4.03'
$ws.Range('E14').Value = '  -0.19%  '
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').Value = '64.80'
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('D17').Value = '26.649.34'
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('D18').Value = '0.0₃0732'
$ws.Range('E18').Value = '  +0.71%  '
$ws.Range('D19').Value = '208.42'
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('D21').Value = '6.96'
$ws.Range('E21').Value = '  +3.06%  '
$ws.Range('E22').Value = '  +0.25%  '
$ws.Range('D23').Value = '2.31'
$ws.Range('E23').Value = '  -2.83%  '
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('D25').Value = '145.57'
$ws.Range('E25').Value = '  -0.76%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').Value = '7.16'
$ws.Range('E27').Value = '  -1.47%  '
$ws.Range('D28').Value = '0.115'
$ws.Range('E28').Value = '  +0.69%  '
$ws.Range('D29').Value = '15.30'
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').Value = '0.0505'
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('D31').Value = '1.15'
$ws.Range('E31').Value = '  +0.08%  '
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('D33').Value = '0.656'
$ws.Range('E33').Value = '  -1.12%  '
$ws.Range('E34').Value = '  +0.33%  '
$ws.Range('D35').Value = '1.283.01'
$ws.Range('E35').Value = '  -2.75%  '
$ws.Range('E36').Value = '  +1.60%  '
$ws.Range('D37').Value = '1.50'
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('E38').Value = '  -0.19%  '
$ws.Range('D39').Value = '0.844'
$ws.Range('E39').Value = '  +1.77%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('E41').Value = '  +1.28%  '
$ws.Range('E42').Value = '  +1.65%  '
$ws.Range('D43').Value = '0.785'
$ws.Range('E43').Value = '  -0.56%  '
$ws.Range('D44').Value = '63.99'
$ws.Range('E44').Value = '  +1.08%  '
$ws.Range('D45').Value = '0.918'
$ws.Range('E45').Value = '  +9.29%  '
$ws.Range('D46').Value = '1.736.85'
$ws.Range('E46').Value = '  +0.60%  '
$ws.Range('D47').Value = '89.79'
$ws.Range('E47').Value = '  -0.28%  '
$ws.Range('E48').Value = '  -0.52%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0105'
$ws.Range('E49').Value = '  -1.18%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '0.102'
$ws.Range('E50').Value = '  +3.74%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.0506'
$ws.Range('E51').Value = '  -1.07%  '

# Restore the default (unstyled) cell style now that the text values are set,
# matching the workbook's original unstyled data cells.
$priceRange.Style = "Normal"
